# Correction in SA algorithm and 746 logs
# Update the "Fitness" column (C) values for run_10 log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C16").Value = 8101
$ws.Range("C17:C24").Value = 7817
$ws.Range("C25:C26").Value = 7815
$ws.Range("C27").Value = 7769
$ws.Range("C28:C33").Value = 7312
$ws.Range("C34:C252").Value = 7310
